$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "agregue cliente consumidor final" -> append new client ids to the
# existing MEJORAR and PREMIUM client lists.

$ws.Range("B3").Value2 = $ws.Range("B3").Value2 + ".20163"
$ws.Range("B2").Value2 = $ws.Range("B2").Value2 + ".3"

$ws.Range("B3").Select()
